$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1780.6
$ws.Range("I19").Value = 623.8333
$ws.Range("J19").Value = 2551.7778
$ws.Range("K19").Value = 623.8333
$ws.Range("L19").Value = 2551.7778
$ws.Range("M19").Value = -448.8333
$ws.Range("N19").Value = -2901.7778
$ws.Range("H40").Value = 2993
$ws.Range("I40").Value = 2993
$ws.Range("K40").Value = 2993
$ws.Range("M40").Value = -2818
$ws.Range("H51").Value = 6374.25
$ws.Range("I51").Value = 6284.857
$ws.Range("K51").Value = 6284.857
$ws.Range("M51").Value = -5800.857
$ws.Range("H113").Value = 14243.777
$ws.Range("I113").Value = 15749.25
$ws.Range("J113").Value = 2200
$ws.Range("K113").Value = 15749.25
$ws.Range("L113").Value = 2200
$ws.Range("M113").Value = -12495.25
$ws.Range("N113").Value = -8708
$ws.Range("H116").Value = 8261.75
$ws.Range("I116").Value = 14510
$ws.Range("J116").Value = 4096.25
$ws.Range("K116").Value = 14510
$ws.Range("L116").Value = 4096.25
$ws.Range("M116").Value = -11068
$ws.Range("N116").Value = -10980.25
$ws.Range("H133").Value = 39518.223
$ws.Range("J133").Value = 39518.223
$ws.Range("L133").Value = 39518.223
$ws.Range("N133").Value = -49638.223
$ws.Range("H137").Value = 30880.176
$ws.Range("J137").Value = 49376.617
$ws.Range("L137").Value = 148129.851
$ws.Range("N137").Value = -153229.851

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1060066.1
$ws.Range("I2").Value = 1110307.4
$ws.Range("J2").Value = 5000
$ws.Range("K2").Value = 1110307.4
$ws.Range("L2").Value = 5000
$ws.Range("M2").Value = -1110194.4
$ws.Range("N2").Value = -5226
$ws.Range("H32").Value = 2883.3247
$ws.Range("I32").Value = 2268.4092
$ws.Range("K32").Value = 2268.4092
$ws.Range("M32").Value = -1981.4092
$ws.Range("H61").Value = 3896.0715
$ws.Range("I61").Value = 1326.2
$ws.Range("J61").Value = 5323.778
$ws.Range("K61").Value = 1326.2
$ws.Range("L61").Value = 5323.778
$ws.Range("M61").Value = -1114.2
$ws.Range("N61").Value = -5747.778
$ws.Range("H104").Value = 32499.715
$ws.Range("J104").Value = 32499.715
$ws.Range("L104").Value = 32499.715
$ws.Range("N104").Value = -39487.715
$ws.Range("H109").Value = 67379.28999999999
$ws.Range("J109").Value = 67379.28999999999
$ws.Range("L109").Value = 67379.28999999999
$ws.Range("N109").Value = -70153.28999999999
$ws.Range("H112").Value = 29994.25
$ws.Range("J112").Value = 29994.25
$ws.Range("L112").Value = 29994.25
$ws.Range("N112").Value = -32948.25
$ws.Range("H116").Value = 1060066.1
$ws.Range("I116").Value = 1110307.4
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 1110307.4
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = -1108013.4
$ws.Range("N116").Value = -9588
$ws.Range("H119").Value = 50698
$ws.Range("J119").Value = 50698
$ws.Range("L119").Value = 50698
$ws.Range("N119").Value = -60374
$ws.Range("H122").Value = 1586.1562
$ws.Range("I122").Value = 1545.5555
$ws.Range("K122").Value = 4636.666499999999
$ws.Range("M122").Value = -2186.666499999999
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H132").Value = 2012
$ws.Range("I132").Value = 1558.36
$ws.Range("K132").Value = 4675.08
$ws.Range("M132").Value = -2145.08
$ws.Range("H136").Value = 3896.0715
$ws.Range("I136").Value = 1326.2
$ws.Range("J136").Value = 5323.778
$ws.Range("K136").Value = 3978.6
$ws.Range("L136").Value = 15971.334
$ws.Range("M136").Value = -1428.6
$ws.Range("N136").Value = -21071.334

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1060066.1
$ws.Range("I3").Value = 1110307.4
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 1110307.4
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = -1110193.4
$ws.Range("N3").Value = -5228
$ws.Range("H81").Value = 38395.2
$ws.Range("J81").Value = 38395.2
$ws.Range("L81").Value = 38395.2
$ws.Range("N81").Value = -40517.2
$ws.Range("H84").Value = 38395.2
$ws.Range("J84").Value = 38395.2
$ws.Range("L84").Value = 115185.6
$ws.Range("N84").Value = -125793.6
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H110").Value = 99991
$ws.Range("J110").Value = 99991
$ws.Range("L110").Value = 99991
$ws.Range("N110").Value = -108171
$ws.Range("H130").Value = 59993.75
$ws.Range("J130").Value = 59993.75
$ws.Range("L130").Value = 59993.75
$ws.Range("N130").Value = -70033.75
$ws.Range("H134").Value = 5344.636
$ws.Range("I134").Value = 5620.6484
$ws.Range("K134").Value = 16861.9452
$ws.Range("M134").Value = -14326.9452

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2095.1428
$ws.Range("I31").Value = 999.2857
$ws.Range("J31").Value = 4286.857
$ws.Range("K31").Value = 999.2857
$ws.Range("L31").Value = 4286.857
$ws.Range("M31").Value = -704.2857
$ws.Range("N31").Value = -4876.857
$ws.Range("H34").Value = 2095.1428
$ws.Range("I34").Value = 999.2857
$ws.Range("J34").Value = 4286.857
$ws.Range("K34").Value = 999.2857
$ws.Range("L34").Value = 4286.857
$ws.Range("M34").Value = -797.2857
$ws.Range("N34").Value = -4690.857
$ws.Range("H43").Value = 24163.25
$ws.Range("J43").Value = 24163.25
$ws.Range("L43").Value = 24163.25
$ws.Range("N43").Value = -24531.25
$ws.Range("H101").Value = 24163.25
$ws.Range("J101").Value = 24163.25
$ws.Range("L101").Value = 24163.25
$ws.Range("N101").Value = -30653.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3967.074
$ws.Range("I70").Value = 4036.3
$ws.Range("K70").Value = 4036.3
$ws.Range("M70").Value = -3766.3
$ws.Range("H73").Value = 3967.074
$ws.Range("I73").Value = 4036.3
$ws.Range("K73").Value = 4036.3
$ws.Range("M73").Value = -3100.3
$ws.Range("H98").Value = 15797
$ws.Range("J98").Value = 15797
$ws.Range("L98").Value = 15797
$ws.Range("N98").Value = -21787
$ws.Range("H117").Value = 40975
$ws.Range("J117").Value = 40975
$ws.Range("L117").Value = 40975
$ws.Range("N117").Value = -47859
$ws.Range("H122").Value = 1949.2084
$ws.Range("I122").Value = 1868.6666
$ws.Range("J122").Value = 2083.4443
$ws.Range("K122").Value = 5605.9998
$ws.Range("L122").Value = 6250.3329
$ws.Range("M122").Value = -3155.9998
$ws.Range("N122").Value = -11150.3329
$ws.Range("H123").Value = 20162.5
$ws.Range("J123").Value = 20162.5
$ws.Range("L123").Value = 20162.5
$ws.Range("N123").Value = -25062.5
$ws.Range("H132").Value = 1752453.1
$ws.Range("I132").Value = 6412467.5
$ws.Range("J132").Value = 4947.75
$ws.Range("K132").Value = 19237402.5
$ws.Range("L132").Value = 14843.25
$ws.Range("M132").Value = -19234872.5
$ws.Range("N132").Value = -19903.25
$ws.Range("H141").Value = 52199
$ws.Range("J141").Value = 52199
$ws.Range("L141").Value = 52199
$ws.Range("N141").Value = -62559

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 9995.666999999999
$ws.Range("J104").Value = 9995.666999999999
$ws.Range("L104").Value = 9995.666999999999
$ws.Range("N104").Value = -16983.667
$ws.Range("H122").Value = 4657.6294
$ws.Range("I122").Value = 3250.3157
$ws.Range("K122").Value = 9750.947100000001
$ws.Range("M122").Value = -7300.947100000001
$ws.Range("H127").Value = 43641.6
$ws.Range("J127").Value = 43641.6
$ws.Range("L127").Value = 43641.6
$ws.Range("N127").Value = -53561.6
$ws.Range("H132").Value = 2695.2886
$ws.Range("I132").Value = 1017.64
$ws.Range("J132").Value = 4248.6665
$ws.Range("K132").Value = 3052.92
$ws.Range("L132").Value = 12745.9995
$ws.Range("M132").Value = -522.9200000000001
$ws.Range("N132").Value = -17805.9995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 45356.71
$ws.Range("I122").Value = 66465.95
$ws.Range("J122").Value = 1027.3
$ws.Range("K122").Value = 199397.85
$ws.Range("L122").Value = 3081.9
$ws.Range("M122").Value = -196947.85
$ws.Range("N122").Value = -7981.9
